$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "65.544.52" that Excel would
# otherwise auto-convert into a number; pre-format the column as Text so
# the written values stay literal strings, matching the source data.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "65.544.52"
$ws.Range("E2").Value = "  +1.25%  "

$ws.Range("D3").Value = "3.371.43"
$ws.Range("E3").Value = "  +0.98%  "

$ws.Range("E4").Value = "  -0.18%  "

$ws.Range("D5").Value = "182.35"
$ws.Range("E5").Value = "  +0.04%  "

$ws.Range("D6").Value = "539.32"
$ws.Range("E6").Value = "  +1.10%  "

$ws.Range("D7").Value = "0.603"
$ws.Range("E7").Value = "  -0.89%  "

$ws.Range("D8").Value = "3.365.30"
$ws.Range("E8").Value = "  +1.00%  "

$ws.Range("E9").Value = "  -0.15%  "

$ws.Range("D10").Value = "0.627"
$ws.Range("E10").Value = "  +1.38%  "

$ws.Range("D11").Value = "55.19"
$ws.Range("E11").Value = "  -6.70%  "

$ws.Range("D12").Value = "0.144"
$ws.Range("E12").Value = "  +6.15%  "

$ws.Range("D13").Value = "0.0000267"
$ws.Range("E13").Value = "  +1.60%  "

$ws.Range("D14").Value = "9.23"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "3.903.65"
$ws.Range("E15").Value = "  +0.78%  "

$ws.Range("E16").Value = "  +1.60%  "

$ws.Range("D17").Value = "3.363.88"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("D18").Value = "18.00"
$ws.Range("E18").Value = "  +1.71%  "

$ws.Range("D19").Value = "65.705.22"
$ws.Range("E19").Value = "  +1.43%  "

$ws.Range("D20").Value = "11.47"
$ws.Range("E20").Value = "  +1.76%  "

$ws.Range("E21").Value = "  +0.91%  "

$ws.Range("D22").Value = "392.59"
$ws.Range("E22").Value = "  +3.84%  "

$ws.Range("D23").Value = "11.78"
$ws.Range("E23").Value = "  +4.18%  "

$ws.Range("D24").Value = "4.26"
$ws.Range("E24").Value = "  +7.51%  "

$ws.Range("D25").Value = "83.35"
$ws.Range("E25").Value = "  +2.45%  "

$ws.Range("D26").Value = "3.79"
$ws.Range("E26").Value = "  -1.39%  "

$ws.Range("B27").Value = "ImmutableX"
$ws.Range("C27").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D27").Value = "2.84"
$ws.Range("E27").Value = "  +4.98%  "

$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").Value = "5.99"
$ws.Range("E28").Value = "  -1.74%  "

$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D30").Value = "8.47"
$ws.Range("E30").Value = "  -0.19%  "

$ws.Range("D31").Value = "29.55"
$ws.Range("E31").Value = "  +1.09%  "

$ws.Range("D32").Value = "663.79"
$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("E33").Value = "  +0.13%  "

$ws.Range("D34").Value = "11.44"
$ws.Range("E34").Value = "  +0.51%  "

$ws.Range("E35").Value = "  +1.02%  "

$ws.Range("D36").Value = "58.01"
$ws.Range("E36").Value = "  -2.74%  "

$ws.Range("D37").Value = "37.67"
$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("D39").Value = "0.396"
$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").Value = "0.0₃0780"
$ws.Range("E40").Value = "  +9.08%  "

$ws.Range("D41").Value = "2.76"
$ws.Range("E41").Value = "  +9.20%  "

$ws.Range("D42").Value = "3.29"
$ws.Range("E42").Value = "  +16.36%  "

$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("D44").Value = "0.129"
$ws.Range("E44").Value = "  +1.51%  "

$ws.Range("D45").Value = "3.013.00"
$ws.Range("E45").Value = "  +2.21%  "

$ws.Range("D46").Value = "2.77"
$ws.Range("E46").Value = "  +1.26%  "

$ws.Range("D47").Value = "0.0412"
$ws.Range("E47").Value = "  +2.17%  "

$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.21"
$ws.Range("E48").Value = "  +4.04%  "

$ws.Range("B49").Value = "WEMIXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D49").Value = "2.72"
$ws.Range("E49").Value = "  +2.04%  "

$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "8.92"
$ws.Range("E50").Value = "  +11.43%  "

$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").Value = "0.128"
$ws.Range("E51").Value = "  +0.50%  "
